$d = $word.ActiveDocument

# 1. Merge the split runs that make up the operation signature into a single run
#    by replacing the (cross-run) found text with itself.
$d.Content.Find.Execute(
    "opretAfskrivning(navn : String, afskrivningsmetode : afskrivningsmetode)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "opretAfskrivning(navn : String, afskrivningsmetode : afskrivningsmetode)", 2) | Out-Null

# 2. Precondition text tweaks
$d.Content.Find.Execute(
    "En høk studerende h eksisterer", $true, $false, $false, $false, $false, $true, 1, $false,
    "En instans h af HØK eksisterer", 2) | Out-Null

$d.Content.Find.Execute(
    "h har ingen associationer", $true, $false, $false, $false, $false, $true, 1, $false,
    "h har ingen associationer til instanser af Afskrivning", 2) | Out-Null

# 3. Postcondition text tweaks
$d.Content.Find.Execute(
    "nyAfskrivning oprettet", $true, $false, $false, $false, $false, $true, 1, $false,
    "En instans nyAfskrivning af Afskrivning blev oprettet", 2) | Out-Null

$d.Content.Find.Execute(
    "nyAfskrivning har navn", $true, $false, $false, $false, $false, $true, 1, $false,
    "nyAfskrivning.navn blev sat til navn", 2) | Out-Null

$d.Content.Find.Execute(
    "nyAfskrivning har afskrivnings metode", $true, $false, $false, $false, $false, $true, 1, $false,
    "nyAfskrivning blev sat til at blive beregnet med afskrivningsmetode", 2) | Out-Null

$d.Content.Find.Execute(
    "nyAfskrivning klar til værdi indtastelse", $true, $false, $false, $false, $false, $true, 1, $false,
    "nyAfskrivning.afskrivningsværdi blev sat til 0", 2) | Out-Null

# 4. Add a new bullet paragraph "h blev sat til at aflæse nyAfskrivning" right after the
#    "afskrivningsværdi blev sat til 0" bullet (paragraph 13 at this point).
$pValue = $d.Paragraphs(13)
$pValue.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs(14)
$pNew.Range.Text = "h blev sat til at aflæse nyAfskrivning"

# 5. Insert two clean blank paragraphs before the final (originally last, still-empty) paragraph,
#    which is now paragraph 15. Splitting its own paragraph mark in place (mark -> mark) grows the
#    body by one clean empty paragraph without leaving stray runs behind.
$pFinal = $d.Paragraphs(15)
$pFinal.Range.Find.Execute("^p", $true, $false, $false, $false, $false, $true, 1, $false, "^p", 2) | Out-Null

$pFinal2 = $d.Paragraphs(16)
$pFinal2.Range.Find.Execute("^p", $true, $false, $false, $false, $false, $true, 1, $false, "^p", 2) | Out-Null

# 6. Re-home the _GoBack bookmark from the precondition bullet onto the new, second-to-last blank
#    paragraph (the one right before the document's trailing empty paragraph).
$pBookmark = $d.Paragraphs(16)
$d.Bookmarks.Add("_GoBack", $pBookmark.Range)
